# Refresh the cryptos price/volume table (Price = column D, Volume(1h) = column E).
# A few D values are single-decimal numeric-looking strings ("593.97", "30.00", ...);
# prefixing them with a literal apostrophe forces Excel to keep them as text
# (matching the original inline-string cells) instead of silently coercing to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.504.03'
$ws.Range("E2").Value = '  +0.72%  '
$ws.Range("D3").Value = '3.761.66'
$ws.Range("E3").Value = '  -0.47%  '
$ws.Range("D5").Value = '''593.97'
$ws.Range("E5").Value = '  -0.53%  '
$ws.Range("D6").Value = '''167.19'
$ws.Range("E6").Value = '  -1.53%  '
$ws.Range("D7").Value = '3.760.28'
$ws.Range("E7").Value = '  -0.48%  '
$ws.Range("E9").Value = '  -1.00%  '
$ws.Range("E10").Value = '  -2.79%  '
$ws.Range("E11").Value = '  -1.62%  '
$ws.Range("E12").Value = '  -1.11%  '
$ws.Range("E13").Value = '  -6.93%  '
$ws.Range("D14").Value = '''36.10'
$ws.Range("E14").Value = '  -1.69%  '
$ws.Range("D15").Value = '4.396.32'
$ws.Range("E15").Value = '  -0.43%  '
$ws.Range("D16").Value = '3.764.64'
$ws.Range("E16").Value = '  -0.56%  '
$ws.Range("D17").Value = '68.520.87'
$ws.Range("E17").Value = '  +0.88%  '
$ws.Range("D18").Value = '''17.94'
$ws.Range("E18").Value = '  -4.80%  '
$ws.Range("E19").Value = '  +0.80%  '
$ws.Range("E20").Value = '  -3.32%  '
$ws.Range("D21").Value = '''10.75'
$ws.Range("E21").Value = '  +1.50%  '
$ws.Range("D22").Value = '''465.51'
$ws.Range("E22").Value = '  -0.54%  '
$ws.Range("E23").Value = '  -3.32%  '
$ws.Range("D24").Value = '''84.08'
$ws.Range("E24").Value = '  +0.38%  '
$ws.Range("E25").Value = '  -2.34%  '
$ws.Range("D26").Value = '''2.19'
$ws.Range("E26").Value = '  -2.76%  '
$ws.Range("E27").Value = '  -1.68%  '
$ws.Range("E28").Value = '  -4.57%  '
$ws.Range("D30").Value = '3.910.73'
$ws.Range("E30").Value = '  -0.45%  '
$ws.Range("E31").Value = '  -4.82%  '
$ws.Range("D32").Value = '''7.33'
$ws.Range("E32").Value = '  -3.62%  '
$ws.Range("D33").Value = '''30.00'
$ws.Range("E33").Value = '  -1.84%  '
$ws.Range("E34").Value = '  -3.15%  '
$ws.Range("D35").Value = '''9.21'
$ws.Range("E35").Value = '  -0.47%  '
$ws.Range("D37").Value = '3.715.77'
$ws.Range("E37").Value = '  -0.58%  '
$ws.Range("E38").Value = '  -3.87%  '
$ws.Range("D39").Value = '''3.40'
$ws.Range("E39").Value = '  -8.82%  '
$ws.Range("E40").Value = '  -0.54%  '
$ws.Range("E41").Value = '  -0.45%  '
$ws.Range("E42").Value = '  -1.02%  '
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").Value = '''44.04'
$ws.Range("E45").Value = '  +8.93%  '
$ws.Range("E46").Value = '  -3.55%  '
$ws.Range("D47").Value = '''46.82'
$ws.Range("E48").Value = '  -1.94%  '
$ws.Range("E49").Value = '  -2.42%  '
$ws.Range("D50").Value = '''145.49'
$ws.Range("E50").Value = '  +2.51%  '
$ws.Range("D51").Value = '''390.17'
$ws.Range("E51").Value = '  -3.09%  '
